$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.771609
$ws.Range("H2").Value = 26.314827
$ws.Range("I2").Value = 0.2200338127677125
$ws.Range("J2").Value = 0.2200338127677125
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 535.4331034879237
$ws.Range("R2").Value = 4818.897931391314
$ws.Range("S2").Value = 0.0449664061553634
$ws.Range("T2").Value = 0.04496640615536341

$ws.Range("G3").Value = 8.771609
$ws.Range("H3").Value = 26.314827
$ws.Range("I3").Value = 0.2200338127677125
$ws.Range("J3").Value = 0.2200338127677125
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 932.5489267957939
$ws.Range("R3").Value = 8392.940341162146
$ws.Range("S3").Value = 0.07831673747641883
$ws.Range("T3").Value = 0.07831673747641886

$ws.Range("G4").Value = 8.771609
$ws.Range("H4").Value = 26.314827
$ws.Range("I4").Value = 0.2200338127677125
$ws.Range("J4").Value = 0.2200338127677125
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 1152.04917335395
$ws.Range("R4").Value = 10368.44256018556
$ws.Range("S4").Value = 0.09675066913593022
$ws.Range("T4").Value = 0.09675066913593025

$ws.Range("I5").Value = 0.583164828467109
$ws.Range("J5").Value = 0.583164828467109
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 1419.080776829431
$ws.Range("R5").Value = 12771.72699146488
$ws.Range("S5").Value = 0.1191763493189023
$ws.Range("T5").Value = 0.1191763493189023

$ws.Range("I6").Value = 0.583164828467109
$ws.Range("J6").Value = 0.583164828467109
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.2075661290510583
$ws.Range("T6").Value = 0.2075661290510583

$ws.Range("I7").Value = 0.583164828467109
$ws.Range("J7").Value = 0.583164828467109
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 3053.324169198851
$ws.Range("R7").Value = 27479.91752278967
$ws.Range("S7").Value = 0.2564223500971484
$ws.Range("T7").Value = 0.2564223500971484

$ws.Range("G8").Value = 7.845451333333334
$ws.Range("H8").Value = 23.536354
$ws.Range("I8").Value = 0.1968013587651783
$ws.Range("J8").Value = 0.1968013587651783
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 478.8989517966585
$ws.Range("R8").Value = 4310.090566169927
$ws.Range("S8").Value = 0.04021859058318765
$ws.Range("T8").Value = 0.04021859058318766

$ws.Range("G9").Value = 7.845451333333334
$ws.Range("H9").Value = 23.536354
$ws.Range("I9").Value = 0.1968013587651783
$ws.Range("J9").Value = 0.1968013587651783
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 834.0849690323214
$ws.Range("R9").Value = 7506.764721290893
$ws.Range("S9").Value = 0.07004759930095913
$ws.Range("T9").Value = 0.07004759930095915

$ws.Range("G10").Value = 7.845451333333334
$ws.Range("H10").Value = 23.536354
$ws.Range("I10").Value = 0.1968013587651783
$ws.Range("J10").Value = 0.1968013587651783
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 1030.409098622079
$ws.Range("R10").Value = 9273.681887598712
$ws.Range("S10").Value = 0.08653516888103152
$ws.Range("T10").Value = 0.08653516888103155
